# Add a "Round" (replicate number) column to the "Average for R" summary
# sheet, ahead of the existing SampleID/dose column, and tweak the
# remembered selections on a couple of sheets (blocking = the 1/2/3
# replicate round each averaged triplicate came from).

$wb = $excel.ActiveWorkbook

$wsRound = $wb.Worksheets.Item("Average for R")
$wsData  = $wb.Worksheets.Item("Sheet2")

# Insert a brand-new column A; everything that used to live in A:K slides
# over to B:L automatically.
$wsRound.Columns.Item(1).Insert()

$wsRound.Range("A1").Value = "Round"

$roundNumbers = @(1, 2, 3, 1, 2, 3, 1, 2, 3)
for ($i = 0; $i -lt $roundNumbers.Length; $i++) {
    $wsRound.Cells.Item($i + 2, 1).Value = $roundNumbers[$i]
}

# Restore the remembered cursor positions that Excel re-saved along with
# the edit (Sheet2 lost its frozen-ish topLeftCell scroll, and both
# sheets' selections moved).
$wsData.Activate()
$wsData.Range("N13").Select()

$wsRound.Activate()
$wsRound.Range("C11").Select()
